# "Button to add note column" - update header fields (name/cell), mark
# "php"/"Jquery" with spell-check proofErr wrapping, and add a _GoBack
# bookmark, per the supplied OOXML diff.

$d = $word.ActiveDocument

function Set-ParagraphRuns($ParagraphIndex, $RunsXml) {
    $full = $d.Paragraphs($ParagraphIndex).Range
    # Exclude the trailing paragraph-mark character so the paragraph's
    # own <w:pPr> (alignment, numbering, paragraph-mark rPr, etc.) is
    # left completely untouched - only the run content is replaced.
    $target = $d.Range($full.Start, $full.End - 1)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg)
}

# 1) "NOMBRE:____..." -> "NOMBRE:" + " Carlos Jaramillo Corrales"
$runs3 = '<w:r><w:t>NOMBRE:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Carlos Jaramillo Corrales</w:t></w:r>'
Set-ParagraphRuns 3 $runs3

# 2) "CELULAR:_________" + "____" + "_________" -> "CELULAR" + ": 317 534 5577"
#    plus a _GoBack bookmark marking the last edit location.
$runs4 = '<w:r><w:t>CELULAR</w:t></w:r>' +
    '<w:r><w:t>: 317 534 5577</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>'
Set-ParagraphRuns 4 $runs4

# 3) "Se tiene un script en php al cual..." -> split out "php" with
#    spell-check proofErr markers around it.
$runs5 = '<w:r w:rsidRPr="001948BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Se tiene un script en </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>php</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> al cual se le deben hacer las siguientes modificaciones:</w:t></w:r>'
Set-ParagraphRuns 5 $runs5

# 4) "...la escala presentada en el script php" -> split out "php" with
#    spell-check proofErr markers, keeping neighbouring runs intact.
$aAcute = [char]0x00E1
$iAcute = [char]0x00ED
$runs8 = '<w:r w:rsidRPr="001948BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Crear </w:t></w:r>' +
    '<w:r w:rsidR="00C6488F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">otra </w:t></w:r>' +
    '<w:r w:rsidRPr="001948BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>columna para mostrar el equivalente del promedio de las notas con las siguientes validaciones: Si la nota es inferior al rango 2 deber' + $aAcute + ' mostrar Bajo, si la nota es inferior a rango 3 deber' + $aAcute + ' mostrar b' + $aAcute + 'sico y as' + $iAcute + ' sucesivamente de acuerdo a</w:t></w:r>' +
    '<w:r w:rsidR="00CC7DD9"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="001948BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">la escala presentada en el script </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>php</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="000A2872" w:rsidRPr="001948BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, se deben tener en cuenta los colores asignados a cada nota. </w:t></w:r>'
Set-ParagraphRuns 8 $runs8

# 5) " de Jquery" -> " de " + "Jquery" wrapped in spell-check proofErr markers.
$oAcute = [char]0x00F3
$runs11 = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Obligatorio u</w:t></w:r>' +
    '<w:r w:rsidR="008E673E"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>sar Ajax</w:t></w:r>' +
    '<w:r w:rsidR="00791D9A"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Jquery</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="008E673E"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> para grabar la informaci' + $oAcute + 'n.</w:t></w:r>'
Set-ParagraphRuns 11 $runs11
